# Update coin Price (D) and Volume(1h) (E) columns for the
# Mon Jan 23 06:42:37 UTC 2023 GitHub Actions symbol-list refresh.
#
# Values must land as literal text (matching the inlineStr cells already
# in the sheet) rather than be auto-coerced to numbers/percentages by
# Excel, and must not pick up a new cell style in the process. Writing
# through a helper cell's text-producing formula and then pasting
# "values only" (PasteSpecial xlPasteValues = -4163) achieves both: the
# pasted result is plain text and the destination cell keeps its original
# (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$updates = @(
    @{ Cell = "D2"; Value = "305.43" },
    @{ Cell = "E2"; Value = "1.23%" },
    @{ Cell = "D3"; Value = "36.05" },
    @{ Cell = "E3"; Value = "-3.98%" },
    @{ Cell = "D4"; Value = "5.100" },
    @{ Cell = "E4"; Value = "1.90%" },
    @{ Cell = "D5"; Value = "0.07867" },
    @{ Cell = "E5"; Value = "0.14%" },
    @{ Cell = "D6"; Value = "2.171" },
    @{ Cell = "E6"; Value = "-2.79%" },
    @{ Cell = "D7"; Value = "7.926" },
    @{ Cell = "E7"; Value = "-1.17%" },
    @{ Cell = "D8"; Value = "0.9185" },
    @{ Cell = "E8"; Value = "0.87%" },
    @{ Cell = "D9"; Value = "0.09696" },
    @{ Cell = "E9"; Value = "4.97%" },
    @{ Cell = "D10"; Value = "0.1865" },
    @{ Cell = "E10"; Value = "-0.98%" },
    @{ Cell = "D11"; Value = "0.08693" },
    @{ Cell = "E11"; Value = "1.98%" },
    @{ Cell = "D12"; Value = "0.03490" },
    @{ Cell = "E12"; Value = "-1.14%" },
    @{ Cell = "D13"; Value = "0.09926" },
    @{ Cell = "E13"; Value = "-0.06%" },
    @{ Cell = "D14"; Value = "0.001447" },
    @{ Cell = "E14"; Value = "-2.67%" },
    @{ Cell = "D15"; Value = "0.005728" },
    @{ Cell = "E15"; Value = "1.12%" },
    @{ Cell = "D16"; Value = "3.460" },
    @{ Cell = "E16"; Value = "-0.31%" },
    @{ Cell = "D17"; Value = "4.100" },
    @{ Cell = "E17"; Value = "1.88%" },
    @{ Cell = "D18"; Value = "2.395" },
    @{ Cell = "E18"; Value = "11.25%" },
    @{ Cell = "D19"; Value = "0.3424" },
    @{ Cell = "E19"; Value = "-1.13%" },
    @{ Cell = "E20"; Value = "-0.50%" },
    @{ Cell = "D21"; Value = "4.854" },
    @{ Cell = "E21"; Value = "1.55%" },
    @{ Cell = "D22"; Value = "0.2201" },
    @{ Cell = "E22"; Value = "0.04%" },
    @{ Cell = "D23"; Value = "0.04557" },
    @{ Cell = "E23"; Value = "-1.92%" },
    @{ Cell = "D24"; Value = "0.005087" },
    @{ Cell = "E24"; Value = "14.31%" },
    @{ Cell = "D25"; Value = "0.001233" },
    @{ Cell = "E25"; Value = "0.44%" },
    @{ Cell = "D26"; Value = "0.0001400" },
    @{ Cell = "E26"; Value = "7.90%" },
    @{ Cell = "D27"; Value = "0.0004751" },
    @{ Cell = "E27"; Value = "0.21%" },
    @{ Cell = "D39"; Value = "0.01839" },
    @{ Cell = "E39"; Value = "3.93%" },
    @{ Cell = "D40"; Value = "0.04775" },
    @{ Cell = "E40"; Value = "1.08%" },
    @{ Cell = "D41"; Value = "0.007718" },
    @{ Cell = "E41"; Value = "-1.48%" },
    @{ Cell = "D42"; Value = "0.1397" },
    @{ Cell = "E42"; Value = "0.33%" },
    @{ Cell = "D43"; Value = "0.007831" },
    @{ Cell = "E43"; Value = "2.26%" },
    @{ Cell = "D44"; Value = "0.002230" },
    @{ Cell = "E44"; Value = "0.63%" },
    @{ Cell = "E45"; Value = "10.66%" },
    @{ Cell = "D46"; Value = "0.00006396" },
    @{ Cell = "E46"; Value = "6.92%" },
    @{ Cell = "E47"; Value = "0.19%" },
    @{ Cell = "D48"; Value = "0.0005801" },
    @{ Cell = "E48"; Value = "0.01%" },
    @{ Cell = "D49"; Value = "24.50" },
    @{ Cell = "E49"; Value = "182.58%" },
    @{ Cell = "E50"; Value = "-25.49%" },
    @{ Cell = "D51"; Value = "0.00002101" },
    @{ Cell = "E51"; Value = "0.19%" }
)

foreach ($u in $updates) {
    $helper.Formula = "=""" + $u.Value + """"
    $helper.Copy() | Out-Null
    $ws.Range($u.Cell).PasteSpecial(-4163)
}

$helper.Clear()
$excel.CutCopyMode = $false
